# edit.ps1
# Appends two new Google-Forms survey response rows (34 and 35) to the
# "Form Responses 1" sheet, matching two new participants ("Cameron" and
# "Folivora20") that were collected after the workbook was last saved.
#
# Strategy:
#  - Column formats for the new rows are inherited from row 5 (an existing
#    data row) via a format-only paste, which reproduces the same per-column
#    cell styles (date format in column A, normal text style elsewhere) that
#    Excel/Sheets applies automatically when a form response is appended.
#  - Cell values are then written on top of the copied formatting.
#  - Two of the answers on row 35 (columns Z and AE) were left blank by the
#    respondent, so those cells are explicitly cleared/left empty.
#  - The bottom-right corner of the data (CC35) is left as the active
#    selection, mirroring where the editor's cursor ended up after typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Data for new rows 34 and 35 (Google-Forms export rows appended to the sheet) ----
$cols = @('A', 'B', 'C', 'D', 'E', 'F', 'G', 'H', 'I', 'J', 'K', 'L', 'M', 'N', 'O', 'P', 'Q', 'R', 'S', 'T', 'U', 'V', 'W', 'X', 'Y', 'Z', 'AA', 'AB', 'AC', 'AD', 'AE', 'AF', 'AG', 'AH', 'AI', 'AJ', 'AK', 'AL', 'AM', 'AN', 'AO', 'AP', 'AQ', 'AR', 'AS', 'AT', 'AU', 'AV', 'AW', 'AX', 'AY', 'AZ', 'BA', 'BB', 'BC', 'BD', 'BE', 'BF', 'BG', 'BH', 'BI', 'BJ', 'BK', 'BL', 'BM', 'BN', 'BO', 'BP', 'BQ', 'BR', 'BS', 'BT', 'BU', 'BV', 'BW', 'BX', 'BY', 'BZ', 'CA', 'CB')

$row34 = @(44208.805103692124, 'Cameron', 'Very infrequently', 'Somewhat frequently', 'Somewhat frequently', 'Somewhat frequently', 'Somewhat frequently', 'Somewhat frequently', 'Somewhat frequently', 'Very frequently', 'Very frequently', 'Somewhat frequently', 'Very infrequently', 'Very frequently', 'Somewhat frequently', 'Somewhat frequently', 'Very frequently', 'Often', 'Often', 'Often', 'Often', 'Often', 'Often', 'Often', 'Often', 'Often', 'Sometimes', 'Sometimes', 'Sometimes', 'Slightly agree', 'Slightly agree', 'Slightly agree', 'Slightly agree', 'Slightly agree', 'Slightly agree', 'Slightly agree', 'Slightly agree', 'No', 'Yes', 'Yes', 'Yes', 'Yes', 'Yes', 'No', 'Yes', 'Yes', 'Yes', 'Yes', 'Yes', 'Yes', 'Yes', 'Yes', 'Yes', 'No', 'Yes', 'Yes', 'No', 'Yes', 'Yes', 'Rather true', 'Rather true', 'Not true', 'Hardly true', 'Hardly true', 'Rather true', 'Rather true', 'Rather true', 'Rather true', 'Rather true', 'Most of the time', 'Most of the time', 'Most of the time', 'Most of the time', 'Some of the time', 'Most of the time', 'A little of the time', 8, 7, 7, 'You were a little better than other workers')
$row35 = @(44209.580218356481, 'Folivora20', 'Somewhat infrequently', 'Very infrequently', 'Somewhat frequently', 'Somewhat frequently', 'Very frequently', 'Somewhat frequently', 'Somewhat infrequently', 'Somewhat frequently', 'Somewhat infrequently', 'Somewhat infrequently', 'Very infrequently', 'Somewhat frequently', 'Somewhat frequently', 'Somewhat frequently', 'Very frequently', 'Often', 'Sometimes', 'Very often or always', 'Sometimes', 'Often', 'Rarely', 'Often', 'Rarely', $null, 'Sometimes', 'Sometimes', 'Often', 'Agree', $null, 'Slightly agree', 'Slightly agree', 'Slightly agree', 'Slightly agree', 'Agree', 'Agree', 'Yes', 'Yes', 'Yes', 'Yes', 'No', 'Yes', 'No', 'No', 'Yes', 'No', 'Yes', 'No', 'Yes', 'Yes', 'No', 'Yes', 'Yes', 'Yes', 'Yes', 'No', 'No', 'Yes', 'Rather true', 'Rather true', 'Hardly true', 'Hardly true', 'Rather true', 'Hardly true', 'Rather true', 'Rather true', 'Hardly true', 'Hardly true', 'Some of the time', 'A little of the time', 'Some of the time', 'A little of the time', 'A little of the time', 'Some of the time', 'Some of the time', 7, 8, 7, 'You were about average')

$newRows = @{ 34 = $row34; 35 = $row35 }

foreach ($r in @(34, 35)) {
    # Inherit the column-by-column cell formatting used by the existing data
    # rows (date format in col A, normal text style elsewhere) before writing
    # the new values.
    $ws.Range("A5:CB5").Copy()
    $ws.Range(("A{0}:CB{0}" -f $r)).PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false

    $data = $newRows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $value = $data[$i]
        $addr = ("{0}{1}" -f $cols[$i], $r)
        if ($null -eq $value) {
            $ws.Range($addr).Value = $null
        } else {
            $ws.Range($addr).Value = $value
        }
    }
}

# Restore the cursor/selection to the last cell of the newly entered data,
# same place it would be after typing in the final value.
$ws.Range("CC35").Select()
